$wb = $excel.ActiveWorkbook

# Add a new worksheet "Resources" at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Resources"

# Populate the new sheet's data (order matches the original authoring order
# so that the shared-string table indices line up with the target workbook)
$newSheet.Range("B1").Value = "Link"
$newSheet.Range("A1").Value = "Title"

$newSheet.Range("A2").Value = "How to use MMC/SDC"
$newSheet.Range("B2").Value = "http://elm-chan.org/docs/mmc/mmc_e.html"

$newSheet.Range("B3").Value = "https://github.com/kiwih/cubeide-sd-card"
$newSheet.Range("A3").Value = "SD Card SPI Firmware"

$newSheet.Range("A4").Value = "MicroSD SPI Breakout Board"
$newSheet.Range("B4").Value = "https://cdn-learn.adafruit.com/downloads/pdf/adafruit-microsd-spi-sdio.pdf"

$newSheet.Range("A5").Value = "SD Card Over SPI Tutorial"
$newSheet.Range("B5").Value = "https://01001000.xyz/2020-08-09-Tutorial-STM32CubeIDE-SD-card/"

# Set column widths to match bestFit sizing (26 chars / ~71.57 chars)
$newSheet.Columns.Item(1).ColumnWidth = 25.1
$newSheet.Columns.Item(2).ColumnWidth = 70.6

# Make the new sheet the active sheet/tab
$newSheet.Activate()
$newSheet.Range("C13").Select()
